$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
  @{ A="ECs"; B="Omg"; C="Rtn4rl1"; D="ECs"; E=3; F=1; G=0.8986206666666666; H=2.695862; I=0.1661804693926261; J=0.1661804693926262; K=1; L=0.3333333333333333; M=0.05112433333333333; N=0.153373; O=0.006671445085820153; P=0.006671445085820154; Q=0.04594138250288889; R=0.413472442526; S=0.001108663875888722; T=0.001108663875888722 },
  @{ A="ECs"; B="Omg"; C="Rtn4rl1"; D="FAPs"; E=3; F=1; G=0.8986206666666666; H=2.695862; I=0.1661804693926261; J=0.1661804693926262; K=3; L=1; M=4.806204333333334; N=14.418613; O=0.6271833037313775; P=0.6271833037313777; Q=4.318954542156222; R=38.870590879406; S=0.1042256158092983; T=0.1042256158092984 },
  @{ A="ECs"; B="Omg"; C="Rtn4rl1"; D="MuSCs"; E=3; F=1; G=0.8986206666666666; H=2.695862; I=0.1661804693926261; J=0.1661804693926262; K=3; L=1; M=2.805828666666667; N=8.417486; O=0.3661452511828023; P=0.3661452511828023; Q=2.521375626992445; R=22.692380642932; S=0.06084618970743907; T=0.0608461897074391 },
  @{ A="FAPs"; B="Omg"; C="Rtn4rl1"; D="ECs"; E=3; F=1; G=2.109481; H=6.328443; I=0.3901029163453022; J=0.3901029163453023; K=1; L=0.3333333333333333; M=0.05112433333333333; N=0.153373; O=0.006671445085820153; P=0.006671445085820154; Q=0.1078458098043333; R=0.970612288239; S=0.002602550184215977; T=0.002602550184215977 },
  @{ A="FAPs"; B="Omg"; C="Rtn4rl1"; D="FAPs"; E=3; F=1; G=2.109481; H=6.328443; I=0.3901029163453022; J=0.3901029163453023; K=3; L=1; M=4.806204333333334; N=14.418613; O=0.6271833037313775; P=0.6271833037313777; Q=10.13859672328434; R=91.247370509559; S=0.2446660358686918; T=0.2446660358686919 },
  @{ A="FAPs"; B="Omg"; C="Rtn4rl1"; D="MuSCs"; E=3; F=1; G=2.109481; H=6.328443; I=0.3901029163453022; J=0.3901029163453023; K=3; L=1; M=2.805828666666667; N=8.417486; O=0.3661452511828023; P=0.3661452511828023; Q=5.918842261588667; R=53.269580354298; S=0.1428343302923944; T=0.1428343302923944 },
  @{ A="MuSCs"; B="Omg"; C="Rtn4rl1"; D="ECs"; E=3; F=1; G=2.399397; H=7.198191; I=0.4437166142620716; J=0.4437166142620716; K=1; L=0.3333333333333333; M=0.05112433333333333; N=0.153373; O=0.006671445085820153; P=0.006671445085820154; Q=0.122667572027; R=1.104008148243; S=0.002960231025715454; T=0.002960231025715454 },
  @{ A="MuSCs"; B="Omg"; C="Rtn4rl1"; D="FAPs"; E=3; F=1; G=2.399397; H=7.198191; I=0.4437166142620716; J=0.4437166142620716; K=3; L=1; M=4.806204333333334; N=14.418613; O=0.6271833037313775; P=0.6271833037313777; Q=11.531992258787; R=103.787930329083; S=0.2782916520533874; T=0.2782916520533874 },
  @{ A="MuSCs"; B="Omg"; C="Rtn4rl1"; D="MuSCs"; E=3; F=1; G=2.399397; H=7.198191; I=0.4437166142620716; J=0.4437166142620716; K=3; L=1; M=2.805828666666667; N=8.417486; O=0.3661452511828023; P=0.3661452511828023; Q=6.732296885314001; R=60.590671967826; S=0.1624647311829688; T=0.1624647311829688 }
)

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")
$startRow = 2
for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $startRow + $i
    $rowData = $rows[$i]
    foreach ($c in $cols) {
        $ws.Range(($c + $r)).Value = $rowData[$c]
    }
}

Write-Output "done"